$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 607
$ws1.Range("F7").Value = 52
$ws1.Range("F12").Value = 101
$ws1.Range("F14").Value = 415
$ws1.Range("F17").Value = 11154
$ws1.Range("F18").Value = 5335

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 607
$ws4.Range("F7").Value = 52
$ws4.Range("F14").Value = 101
$ws4.Range("F16").Value = 415
$ws4.Range("F19").Value = 11154
$ws4.Range("F21").Value = 5335
